$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The H column (roboticRNAPrep) cells currently hold Boolean FALSE values
# displayed via a custom "TRUE/FALSE" number format. Convert them to the
# literal text "False" (shared string), matching the rest of the sheet's
# plain-text style, and restyle the column as Text.
$rng = $ws.Range("H2:H27")

# Build the literal text via a formula first (direct Value/Value2 writes of
# "False" get auto-typed back to a Boolean by Excel's input parser), then
# freeze the formula results down to plain values with Copy/PasteSpecial.
$rng.Formula = "=T(""False"")"
$rng.Copy()
$rng.PasteSpecial(-4163)
$rng.NumberFormat = "@"

$ws.Range("H3:H27").Select()
